$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2 value
$ws.Range("B2").Value = 0.5498201444332311

# Add new rows 3-5
$ws.Range("A3").Value = 20
$ws.Range("B3").Value = 0.5908374480150782

$ws.Range("A4").Value = 30
$ws.Range("B4").Value = 0.620589314846667

$ws.Range("A5").Value = 40
$ws.Range("B5").Value = 0.6308227565762184
